# Automatische test-sync: 2025-08-18 20:33:50
#
# Adds a new "Interne taak" log row (row 5) to the Logs sheet, bumps the
# matching "Overig" tally row (row 4) on the Dashboard sheet, and extends
# the conditional formatting + chart series ranges so they keep covering
# the newly added rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append the new row of data (row 5)
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A5").Value = "Interne taak"
$logs.Range("B5").Value = "kwaliteit@testbedrijf123.nl"
$logs.Range("C5").Value = "Leg dit even neer bij Koen."
$logs.Range("D5").Value = "Overig"
$logs.Range("F5").Value = "2025-08-18 20:32:52"
$logs.Range("G5").Value = "Nee"
$logs.Range("H5").Value = "Ja"
$logs.Range("I5").Value = "Nee"
$logs.Range("J5").Value = "Nee"

# Extend the conditional formatting ranges (previously *2:*4) down to row 5
# so the new row is covered the same way the other data rows are.
$logs.Range("D2:D4").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D5"))
$logs.Range("G2:G4").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G5"))
$logs.Range("H2:H4").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H5"))
$logs.Range("I2:I4").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I5"))
$logs.Range("J2:J4").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J5"))

# ---------------------------------------------------------------------
# 2. Dashboard sheet: add the "Overig" tally row (row 4)
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A4").Value = "Overig"
$dash.Range("B4").Value = 1

# ---------------------------------------------------------------------
# 3. Chart: extend the category/value series references to include the
#    freshly added Dashboard row (A2:A3 -> A2:A4, B2:B3 -> B2:B4)
# ---------------------------------------------------------------------
$chartObj = $dash.ChartObjects(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection(1)
$series.Formula = "=SERIES('Dashboard'!`$B`$1,'Dashboard'!`$A`$2:`$A`$4,'Dashboard'!`$B`$2:`$B`$4,1)"
